$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new item "TSL2651T" ---
$ws.Range("A11").Value = "TSL2651T"
$ws.Range("F11").Value = "https://uk.rs-online.com/web/p/ambient-light-sensors/6424367/"

# --- Row 12: new item "mount adaptors" (digikey link string must be registered before the
#     "mount adaptors" text so the shared-string table order matches) ---
$ws.Range("F12").Value = "https://www.digikey.co.uk/product-detail/en/sparkfun-electronics/BOB-00717/1568-1098-ND/5318740?utm_adgroup=&mkwid=ss3oA964O&pcrid=337808698799&pkw=&pmt=&pdv=c&productid=5318740&slid=&gclid=CjwKCAjw-ZvlBRBbEiwANw9UWrJUGs_IqFlKc6eegGW9viSX1XXCmb91nn1dUuqMT3nWPoyzsXy8kBoCCQAQAvD_BwE"
$ws.Range("A12").Value = "mount adaptors"

# --- Row 13: new item "IR LEDs" ---
$ws.Range("A13").Value = "IR LEDs"
$ws.Range("F13").Value = "https://uk.rs-online.com/web/p/ir-leds/6997663/"

# --- Numeric cost / quantity / line-total formulas for the three new rows ---
$ws.Range("B11").Value = 2.646
$ws.Range("C11").Value = 5
$ws.Range("D11").Formula = "=B11*C11"

$ws.Range("B12").Value = 0.73
$ws.Range("C12").Value = 5
$ws.Range("D12").Formula = "=B12*C12"

$ws.Range("B13").Value = 0.364
$ws.Range("C13").Value = 10
$ws.Range("D13").Formula = "=B13*C13"

# --- The old "Total"/"Remaining Budget" row (row 11) moves down to row 14 now that
#     three extra purchased items were inserted above it ---
$ws.Range("E11").ClearContents()
$ws.Range("C14").Value = "Total"
$ws.Range("D14").Formula = "=SUM(D4:D13)"
$ws.Range("E14").Formula = "=75-D14"

# --- Selection / view bookkeeping to mirror the saved workbook state ---
$ws.Activate() | Out-Null
$ws.Range("E22").Select() | Out-Null
